$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    23 = @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45089, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 60, 11000, 11000, 11000, '$/caja 18 kilos empedrada', 'Región del Maule', 611, 18)
    24 = @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45089, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 50, 9000, 9000, 9000, '$/caja 18 kilos empedrada', 'Región del Maule', 500, 18)
    25 = @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45089, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Segunda', 30, 7000, 7000, 7000, '$/caja 18 kilos empedrada', 'Región del Maule', 389, 18)
    26 = @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45070, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 60, 10000, 10000, 10000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 556, 18)
    27 = @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45062, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 50, 13000, 13000, 13000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 722, 18)
    28 = @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45062, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 50, 12000, 12000, 12000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 667, 18)
    29 = @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45085, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 50, 10000, 10000, 10000, '$/caja 18 kilos empedrada', 'Región del Maule', 556, 18)
    30 = @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45033, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Especial', 60, 13000, 13000, 13000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 722, 18)
    31 = @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45033, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 80, 12000, 12000, 12000, '$/caja 18 kilos empedrada', 'Región de O''Higgins', 667, 18)
    32 = @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45076, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 30, 12000, 12000, 12000, '$/caja 15 kilos granel', 'Región de O''Higgins', 800, 15)
    33 = @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45076, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Segunda', 30, 10000, 10000, 10000, '$/caja 15 kilos granel', 'Región de O''Higgins', 667, 15)
    34 = @(7, 'Terminal Hortofrutícola Agro Chillán', 'Ñuble', 45021, 16, 'Fruta', 100104, 'Frutos de pepita', 100104003, 'Membrillo', 'Champion', 'Primera', 50, 12000, 12000, 12000, '$/caja 18 kilos granel', 'Región de O''Higgins', 667, 18)
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    for ($col = 1; $col -le $row.Length; $col++) {
        $ws.Cells.Item([int]$r, $col).Value = $row[$col - 1]
    }
    $ws.Cells.Item([int]$r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Host "Applied weekly Fruta/hortaliza update"
